# Risolto bug HMR e aggiunti altri log per il debug
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Fallimenti per Fragilità" counts for the "hook" rows in both
# groups (LLM: row 6, Analitica: row 13). The dependent percentage (column G,
# shared formula) and the group totals (rows 20/21, columns B & D) recalc
# automatically.
$ws.Range("E6").Value = 3
$ws.Range("E13").Value = 3

# Record the current selection (matches the saved workbook view state).
$ws.Range("E10").Select()
